$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: new entry (2016/12/05 09:41) with custom date+time format
$ws.Range("A4").Value = 42709.40347222222
$ws.Range("A4").NumberFormat = "yyyy/m/d h:mm;@"
$ws.Range("B4").Value = 49
$ws.Range("C4").Value = 80
$ws.Range("E4").Value = 38
$ws.Range("F4").Value = 64
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 6
$ws.Range("M4").Value = 7
$ws.Range("N4").Value = 10

# Row 5: new entry (2016/12/05 14:03) with default short date+time format
$ws.Range("A5").Value = 42709.585416666669
$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
$ws.Range("B5").Value = 56
$ws.Range("C5").Value = 80
$ws.Range("E5").Value = 41
$ws.Range("F5").Value = 56
$ws.Range("I5").Value = 5
$ws.Range("J5").Value = 6
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = 18

# Touch the Q-column totals so the engine recalculates them: Q only
# depends on other formula cells (H, L, P), and the incremental recalc
# doesn't always chain through a second level of shared formulas.
$ws.Range("Q4").Formula = "=SUM(H4,L4,P4)"
$ws.Range("Q5").Formula = "=SUM(H5,L5,P5)"

# Move the active selection to C15, matching the saved view state
$ws.Range("C15").Select()
